$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B10").Value = 59408
$ws.Range("C10").Value = "SIG-3W Lilliput LED Torch &amp; Table Lamp"
$ws.Range("D10").Value = 388.17
$ws.Range("E10").Value = 463.78
$ws.Range("F10").Value = 9
$ws.Range("G10").Value = 3493.53
$ws.Range("B11").Value = 47438
$ws.Range("C11").Value = "SIG-3w Lilliput LED Torch &amp; Table Lamp"
$ws.Range("D11").Value = 401.81
$ws.Range("E11").Value = 480.05
$ws.Range("F11").Value = 2
$ws.Range("G11").Value = 803.62
$ws.Range("F45").Value = 43
$ws.Range("G45").Value = 1431.9
$ws.Range("B46").Value = 26483.77
$ws.Range("F51").Value = 211
$ws.Range("G51").Value = 14778.44
$ws.Range("F54").Value = 47
$ws.Range("G54").Value = 1235.63
$ws.Range("F61").Value = 46
$ws.Range("G61").Value = 1163.34
$ws.Range("F66").Value = 53
$ws.Range("G66").Value = 976.26
$ws.Range("F69").Value = 34
$ws.Range("G69").Value = 2006
$ws.Range("B85").Value = 154363.11
$ws.Range("F95").Value = 8
$ws.Range("G95").Value = 2010
$ws.Range("B103").Value = 24412.29
$ws.Range("F138").Value = 36
$ws.Range("G138").Value = 1897.56
$ws.Range("F141").Value = 47
$ws.Range("G141").Value = 2719.42
$ws.Range("F142").Value = 94
$ws.Range("G142").Value = 13294.42
$ws.Range("F143").Value = 39
$ws.Range("G143").Value = 1554.54
$ws.Range("F146").Value = 33
$ws.Range("G146").Value = 7255.71
$ws.Range("B159").Value = 70101.17
$ws.Range("F171").Value = 115
$ws.Range("G171").Value = 4866.8
$ws.Range("F175").Value = 177
$ws.Range("G175").Value = 5403.81
$ws.Range("F179").Value = 70
$ws.Range("G179").Value = 3116.4
$ws.Range("B180").Value = 35706.93
$ws.Range("F188").Value = 27
$ws.Range("G188").Value = 2508.03
$ws.Range("F190").Value = 41
$ws.Range("G190").Value = 5187.73
$ws.Range("B198").Value = 43243.9
$ws.Range("F221").Value = 113
$ws.Range("G221").Value = 2865.68
$ws.Range("F225").Value = 47
$ws.Range("G225").Value = 3678.22
$ws.Range("F227").Value = 39
$ws.Range("G227").Value = 1485.12
$ws.Range("B228").Value = 34707.97
$ws.Range("F248").Value = 49
$ws.Range("G248").Value = 1165.22
$ws.Range("F261").Value = 95
$ws.Range("G261").Value = 1775.55
$ws.Range("B267").Value = 20762.94
$ws.Range("F284").Value = 18
$ws.Range("G284").Value = 4787.28
$ws.Range("F286").Value = 30
$ws.Range("G286").Value = 5994.3
$ws.Range("F287").Value = 14
$ws.Range("G287").Value = 3230.64
$ws.Range("B292").Value = 41864
$ws.Range("C292").Value = "HAM-THERMOSTEEL 1000 ML WITH PLAIN LID"
$ws.Range("F292").Value = 1
$ws.Range("G292").Value = 672.04
$ws.Range("B293").Value = 56449
$ws.Range("C293").Value = "HAM-Thermosteel 1000 Ml With Plain Lid"
$ws.Range("F293").Value = 25
$ws.Range("G293").Value = 16801
$ws.Range("B298").Value = 107502.19
$ws.Range("F324").Value = 42
$ws.Range("G324").Value = 3289.86
$ws.Range("B349").Value = 140122.1
$ws.Range("F382").Value = 26
$ws.Range("G382").Value = 2855.32
$ws.Range("F409").Value = 158
$ws.Range("G409").Value = 27070.14
$ws.Range("F422").Value = 11
$ws.Range("G422").Value = 1580.26
$ws.Range("B423").Value = 140054.07
$ws.Range("F428").Value = 19
$ws.Range("G428").Value = 2914.79
$ws.Range("B437").Value = 19762.37
$ws.Range("F458").Value = 14
$ws.Range("G458").Value = 1677.62
$ws.Range("B481").Value = 42587.51
$ws.Range("F486").Value = 91
$ws.Range("G486").Value = 5523.7
$ws.Range("F487").Value = 5
$ws.Range("G487").Value = 276.65
$ws.Range("B497").Value = 36041.26
$ws.Range("F526").Value = 751
$ws.Range("G526").Value = 72546.60000000001
$ws.Range("F527").Value = 168
$ws.Range("G527").Value = 6251.28
$ws.Range("F529").Value = 158
$ws.Range("G529").Value = 4248.62
$ws.Range("B532").Value = 151126.38
$ws.Range("F542").Value = 26
$ws.Range("G542").Value = 1134.9
$ws.Range("F547").Value = 12
$ws.Range("G547").Value = 863.64
$ws.Range("F554").Value = 0
$ws.Range("G554").Value = 0
$ws.Range("B556").Value = 14220.61
$ws.Range("F563").Value = 208
$ws.Range("G563").Value = 3352.96
$ws.Range("F565").Value = 135
$ws.Range("G565").Value = 4691.25
$ws.Range("F566").Value = 152
$ws.Range("G566").Value = 2919.92
$ws.Range("B567").Value = 50600.44
$ws.Range("F592").Value = 84
$ws.Range("G592").Value = 2629.2
$ws.Range("F595").Value = 64
$ws.Range("G595").Value = 3205.12
$ws.Range("F609").Value = 31
$ws.Range("G609").Value = 3045.13
$ws.Range("B610").Value = 55201.68
$ws.Range("F621").Value = 233
$ws.Range("G621").Value = 14147.76
$ws.Range("F623").Value = 69
$ws.Range("G623").Value = 5917.44
$ws.Range("F625").Value = 21
$ws.Range("G625").Value = 1350.72
$ws.Range("B638").Value = 145771.25
$ws.Range("F673").Value = 142
$ws.Range("G673").Value = 3703.36
$ws.Range("F674").Value = 95
$ws.Range("G674").Value = 4957.1
$ws.Range("B688").Value = 86237.48
$ws.Range("F713").Value = 23
$ws.Range("G713").Value = 3002.65
$ws.Range("F714").Value = 65
$ws.Range("G714").Value = 11571.3
$ws.Range("F718").Value = 110
$ws.Range("G718").Value = 2992
$ws.Range("B720").Value = 27847.1
$ws.Range("F731").Value = 9
$ws.Range("G731").Value = 2797.92
$ws.Range("B742").Value = 4570.42
$ws.Range("F778").Value = 8
$ws.Range("G778").Value = 572.48
$ws.Range("F780").Value = 100
$ws.Range("G780").Value = 8550
$ws.Range("F781").Value = 7
$ws.Range("G781").Value = 426.51
$ws.Range("F784").Value = 67
$ws.Range("G784").Value = 4794.52
$ws.Range("B785").Value = 14343.51
$ws.Range("F814").Value = 85
$ws.Range("G814").Value = 11928.05
$ws.Range("B815").Value = 37665.98
$ws.Range("F820").Value = 52
$ws.Range("G820").Value = 4241.12
$ws.Range("F821").Value = 118
$ws.Range("G821").Value = 15705.8
$ws.Range("F829").Value = 109
$ws.Range("G829").Value = 5856.57
$ws.Range("B837").Value = 190587.88
$ws.Range("F843").Value = 64
$ws.Range("G843").Value = 6963.84
$ws.Range("F846").Value = 41
$ws.Range("G846").Value = 3462.04
$ws.Range("F861").Value = 297
$ws.Range("G861").Value = 10938.51
$ws.Range("F866").Value = 46
$ws.Range("G866").Value = 2621.54
$ws.Range("B867").Value = 199931.73
$ws.Range("F891").Value = 3
$ws.Range("G891").Value = 1720.68
$ws.Range("B904").Value = 37896.49
$ws.Range("B923").Value = 2549119.61
$ws.Range("B924").Value = 2549119.61
